# "Add files via upload" - refresh of the daily tracker: a new day
# ("29-sep") is appended as the next column on Sheet1, and Sheet3's
# VLOOKUP source table (A20:B36) is refreshed with that day's figures,
# which ripples into Sheet3's VLOOKUP column C2:C18.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# New column CE on Sheet1: header "29-sep" plus per-row values
$ws1.Cells.Item(1, 83).Value = "29-sep"
$ws1.Cells.Item(2, 83).Value = 0
$ws1.Cells.Item(3, 83).Value = 12.157100117135014
$ws1.Cells.Item(4, 83).Value = 9.3589905312070112
$ws1.Cells.Item(5, 83).Value = 18.174507366868241
$ws1.Cells.Item(6, 83).Value = 0
$ws1.Cells.Item(7, 83).Value = 0.051298452948722956
$ws1.Cells.Item(8, 83).Value = 11.733316302379142
$ws1.Cells.Item(9, 83).Value = 13.622350936712492
$ws1.Cells.Item(10, 83).Value = 19.190135605703976
$ws1.Cells.Item(11, 83).Value = 8.7832946709545148
$ws1.Cells.Item(12, 83).Value = 0
$ws1.Cells.Item(13, 83).Value = 3.865789780510545
$ws1.Cells.Item(14, 83).Value = 0
$ws1.Cells.Item(15, 83).Value = 0
$ws1.Cells.Item(16, 83).Value = 9.596215839287412
$ws1.Cells.Item(17, 83).Value = 0
$ws1.Cells.Item(18, 83).Value = 0

# Updated lookup table B20:B36 on Sheet3 (flows into VLOOKUP column C2:C18)
$ws3.Cells.Item(20, 2).Value = 9.6085615041690247
$ws3.Cells.Item(21, 2).Value = 3.865789780510545
$ws3.Cells.Item(22, 2).Value = 11.733316302379142
$ws3.Cells.Item(23, 2).Value = 0.73568217001000313
$ws3.Cells.Item(24, 2).Value = 0.051298452948722956
$ws3.Cells.Item(25, 2).Value = 9.596215839287412
$ws3.Cells.Item(26, 2).Value = 9.3589905312070112
$ws3.Cells.Item(27, 2).Value = 18.174507366868241
$ws3.Cells.Item(28, 2).Value = 4.7015145853429274
$ws3.Cells.Item(29, 2).Value = 12.157100117135014
$ws3.Cells.Item(30, 2).Value = 18.437108213476126
$ws3.Cells.Item(31, 2).Value = 5.9943798296302946
$ws3.Cells.Item(32, 2).Value = 2.343282153753885
$ws3.Cells.Item(33, 2).Value = 19.190135605703976
$ws3.Cells.Item(34, 2).Value = 8.7832946709545148
$ws3.Cells.Item(35, 2).Value = 13.622350936712492
$ws3.Cells.Item(36, 2).Value = 50.690872523403456

# Selections
$ws3.Activate()
$ws3.Range("C2:C18").Select()
$ws1.Activate()
$ws1.Range("CF6").Select()
